$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes existing rows 10-24 down to 11-25)
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with a new weekly record.
# Same as the record that was previously in row 10, except with a new date.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44494
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Madrigal"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("N10").Value = "$/caja 40 unidades"
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 288
$ws.Range("Q10").Value = 40
$ws.Range("R10").Value = "Hortaliza"
